$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "https://go.mwater.co/cap_leogane"
$ws.Range("D2").Value = "Commune Action Plan: https://go.mwater.co/cap_leogane"

$ws.Range("B3").Value = "---"
$ws.Range("D3").Value = "Project performance: ---"

$ws.Range("B4").Value = "---"
$ws.Range("D4").Value = "Service providers performance: ---"

$ws.Range("B5").Value = "---"
$ws.Range("D5").Value = "Investment status: ---"

$ws.Range("B6").Value = "---"
$ws.Range("D6").Value = "Lessons Learned: ---"
